# edit.ps1 - Recreate the Mau 6 template restructuring described in the commit:
#  - Unmerge cells before writing (avoid merge conflicts)
#  - Restore/expand the template from 11 to 16 columns (A:P), splitting a few
#    combined header fields into their own columns
#  - Remove the stray duplicated header row, fixing the row/STT numbering
#  - Re-merge the banner/footer rows across the new column range

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Unmerge everything first so later writes/merges don't conflict ---
$ws.Cells.UnMerge()

# --- 2. Remove the duplicated header row (old row 17) ---
$ws.Rows("17").Delete()

# --- 3. Extend formatting (style only) from column K into the five new
#        columns L:P for every row, so the new columns inherit the same
#        per-row style used across the rest of that row. ---
$ws.Range("K1:K32").Copy()
$ws.Range("L1:P32").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 4. Rewrite the header row (row 4) with the new, expanded set of
#        column headings in their final order. ---
$ws.Range("A4").Value = "STT"
$ws.Range("B4").Value = "Biển số"
$ws.Range("C4").Value = "Màu biển"
$ws.Range("D4").Value = "Loại xe"
$ws.Range("E4").Value = "Chủ xe"
$ws.Range("F4").Value = "Địa chỉ thường trú, địa chỉ hiện tại của chủ xe"
$ws.Range("G4").Value = "Số khung của xe"
$ws.Range("H4").Value = "Số máy của xe"
$ws.Range("I4").Value = "Số Điện thoại của chủ xe"
$ws.Range("J4").Value = "Số CCCD/mã số thuế của chủ xe"
$ws.Range("K4").Value = "Ngày cấp CCCD"
$ws.Range("L4").Value = "Số GPLX của chủ xe"
$ws.Range("M4").Value = "Ngày cấp GPLX"
$ws.Range("N4").Value = "Cơ quan cấp GPLX"
$ws.Range("O4").Value = "Tình trạng phương tiện (tốt/hỏng)"
$ws.Range("P4").Value = "Ghi chú"

# --- 5. Column widths (approximate to the engine's width rounding) ---
$ws.Columns("A").ColumnWidth = 3.71
$ws.Columns("B").ColumnWidth = 3.71
$ws.Columns("C").ColumnWidth = 9.26
$ws.Columns("D").ColumnWidth = 6.53
$ws.Columns("E").ColumnWidth = 16.8
$ws.Columns("F").ColumnWidth = 25.62
$ws.Columns("G").ColumnWidth = 25.62
$ws.Columns("H").ColumnWidth = 9.98
$ws.Columns("I").ColumnWidth = 9.07
$ws.Columns("J").ColumnWidth = 11.53
$ws.Columns("K").ColumnWidth = 11.53
$ws.Columns("L").ColumnWidth = 11.53
$ws.Columns("M").ColumnWidth = 11.53
$ws.Columns("N").ColumnWidth = 9.53
$ws.Columns("O").ColumnWidth = 9.53
$ws.Columns("P").ColumnWidth = 8.07

# --- 6. Row heights ---
$ws.Rows("1").RowHeight = 58.25
$ws.Rows("2").RowHeight = 13.75
$ws.Rows("3").RowHeight = 23.4
$ws.Rows("4").RowHeight = 111.65
$ws.Rows("5:32").RowHeight = 28

# --- 7. Re-merge the banner/footer rows across the full new column range ---
$ws.Range("A1:P1").Merge()
$ws.Range("A2:P2").Merge()
$ws.Range("A3:P3").Merge()
$ws.Range("A32:P32").Merge()
